$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

# Date header
Replace-Text "2024-05-26 Sunday" "2024-05-27 Monday"

# Table cells (ordered to avoid collisions: 748÷8= -> 222÷8= happens
# before 607÷5= -> 748÷8= is introduced)
Replace-Text "826÷3=" "344÷4="
Replace-Text "580÷6=" "153÷3="
Replace-Text "246÷4=" "858÷7="
Replace-Text "748÷8=" "222÷8="
Replace-Text "716÷5=" "759÷2="
Replace-Text "767÷8=" "162÷6="
Replace-Text "632÷4=" "225÷9="
Replace-Text "370÷9=" "534÷8="
Replace-Text "917÷8=" "675÷7="
Replace-Text "590÷9=" "173÷2="
Replace-Text "838÷8=" "525÷7="
Replace-Text "167÷4=" "525÷8="
Replace-Text "982÷6=" "910÷3="
Replace-Text "707÷6=" "998÷7="
Replace-Text "439÷9=" "562÷2="
Replace-Text "217÷4=" "568÷2="
Replace-Text "221÷8=" "343÷6="
Replace-Text "447÷7=" "251÷6="
Replace-Text "157÷7=" "298÷5="
Replace-Text "607÷7=" "736÷9="
Replace-Text "323÷8=" "945÷7="
Replace-Text "498÷3=" "495÷8="
Replace-Text "607÷5=" "748÷8="
Replace-Text "443÷9=" "702÷9="
Replace-Text "321÷2=" "784÷6="
